$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark (it originally sat right after
# "MP73010" in the title line, a stale artifact of a previous edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# "Ben changing things up!" -> "Krishna changing things up!"
$found = $d.Content.Find.Execute("Ben changing things up!", $true, $false, $false, $false, $false, $true, 1, $false, "Krishna changing things up!", 2)

# Re-create "_GoBack" as a collapsed bookmark right after "Krishna", which is
# where Word leaves it after the most recent text edit. Locating it via a
# fresh Find (rather than a hard-coded offset) keeps this robust regardless
# of which paragraph the sentence lives in.
$locate = $d.Content
$locate.Find.Execute("Krishna changing things up!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterKrishna = $locate.Start + 7
$rng = $d.Range($afterKrishna, $afterKrishna)
$d.Bookmarks.Add("_GoBack", $rng)
